$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sam" (sheet1): fill in "Where?" for week 1 rows, then add a whole
# second week table (rows 9-17) by duplicating the week-1 block's layout and
# formats, same as the user would by copying the block and editing values.
# ---------------------------------------------------------------------------
$sam = $wb.Worksheets.Item("Sam")

# Week 1 "Where?" column entries that were filled in.
$sam.Range("D2").Value = "Rayzor 2nd Floor"
$sam.Range("D3").Value = "My apartment"
$sam.Range("D4").Value = "Rayzor 2nd Floor"

# Build the Week 2 table underneath (row 9 blank spacer, row 10 header,
# rows 11-16 entries, row 17 total) re-using the exact row formatting from
# week 1's table.
$sam.Range("A1:E1").Copy($sam.Range("A10:E10")) | Out-Null
$sam.Range("A2:E2").Copy($sam.Range("A11:E11")) | Out-Null
$sam.Range("A5:E5").Copy($sam.Range("A12:E12")) | Out-Null
$sam.Range("A5:E5").Copy($sam.Range("A13:E13")) | Out-Null
$sam.Range("A5:E5").Copy($sam.Range("A14:E14")) | Out-Null
$sam.Range("A5:E5").Copy($sam.Range("A15:E15")) | Out-Null
$sam.Range("A7:E7").Copy($sam.Range("A16:E16")) | Out-Null
$sam.Range("A8:D8").Copy($sam.Range("A17:D17")) | Out-Null

# New spacer row 9 height matches the rest of the new block.
$sam.Rows.Item(9).RowHeight = 15.75

# Week 2, single entry on row 11.
$sam.Cells.Item(11, 1).Value = 42986.583333333336
$sam.Cells.Item(11, 2).Value = 42986.611111111109
$sam.Cells.Item(11, 3).Formula = "=B11-A11"
$sam.Cells.Item(11, 4).Value = "Rayzor Digital Lab"
$sam.Cells.Item(11, 5).Value = ""

# Shared "=B-A" formula across the blank rows 12:16 (mirrors C3:C7 above).
$sam.Range("C12:C16").Formula = "=B12-A12"

# Week 2 total row.
$sam.Cells.Item(17, 1).Value = "Week 2 Total"
$sam.Cells.Item(17, 3).Formula = "=SUM(C11:C16)"

# Row heights for the new block (header taller to match the wrapped text).
$sam.Rows.Item(10).RowHeight = 27
$sam.Range("A11:E16").RowHeight = 15.75
$sam.Rows.Item(17).RowHeight = 15.75

# Selection / active sheet: "Sam" becomes the active tab, cursor on F7.
$sam.Activate() | Out-Null
$sam.Range("F7").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Travis" (sheet3) is no longer the active tab; update its lingering
# selection to match.
# ---------------------------------------------------------------------------
$travis = $wb.Worksheets.Item("Travis")
$travis.Range("E21").Select() | Out-Null

# Re-activate Sam last so it is the tab shown/selected when the file is saved.
$sam.Activate() | Out-Null
